$d = $word.ActiveDocument

# 1. Remove the "Word version of this document" list item (paragraph + hyperlink)
#    under "Additional resources". A PDF version replaces it (per commit message)
#    but no new paragraph is added in this diff -- the bullet is simply deleted.
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -like "*Word version of this document*") {
        $p.Range.Delete()
    }
}

# 2. Fix wording in the Pitfalls section about inflection points.
$d.Content.Find.Execute(
    "about inflection points. There is any particular relevance to statistics.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "about inflection points, not about statistics.", 2)
